# Insert a new column before column U (21st column) so that existing
# columns U:V shift right to become V:W. This preserves their contents
# and formatting (including the header style) while making room for
# the new "MgCa Coretop modelled temperature" column.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns.Item(21).Insert()

# New header for the inserted column U1.
$ws.Range("U1").Value() = "MgCa Coretop modelled temperature"

# Updated ERSST-based anomaly values in row 2.
$ws.Range("R2").Value() = 27.68
$ws.Range("S2").Value() = -1.182914225260415
$ws.Range("T2").Value() = -0.8977142252604153

# New coretop-modelled-temperature value for the inserted column.
$ws.Range("U2").Value() = 28.2441

# Re-assert the shifted coretop anomaly values (now in V2/W2) so they
# keep their original literal values.
$ws.Range("V2").Value() = -1.7441
$ws.Range("W2").Value() = -1.4589
